# Fix multiplanet coloring translations:
# The "PlanetColor" sheet's Planet 1 / Planet 2 labels incorrectly carried a
# trailing colon (and the French column had been accidentally sharing the
# "Planète 1:"/"Planète 2:" strings with the unrelated PlanetOpacity sheet).
# Strip the trailing colon from every language's translation for the two
# "Planet 1" / "Planet 2" cells so PlanetColor gets its own distinct strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanetColor")

# Row 4 - "Planet 1" label (PLANET_ONE)
$ws.Range("B4").Value = "Planet 1"
$ws.Range("C4").Value = "행성 1"
$ws.Range("D4").Value = "Planeta 1"
$ws.Range("F4").Value = "Planète 1"
$ws.Range("G4").Value = "Hành tinh thứ nhất"
$ws.Range("I4").Value = "Planet 1"

# Row 5 - "Planet 2" label (PLANET_TWO)
$ws.Range("B5").Value = "Planet 2"
$ws.Range("C5").Value = "행성 2"
$ws.Range("D5").Value = "Planeta 2"
$ws.Range("F5").Value = "Planète 2"
$ws.Range("G5").Value = "Hành tinh thứ hai"
$ws.Range("I5").Value = "Planet 2"
